$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Insert three new paragraphs ("One", " Two", " Three") before the
# existing first paragraph ("a").
[void]$tr.InsertBefore("One`r Two`r Three`r")

# Scope to each newly-inserted paragraph's text precisely.
$r1 = $tr.Characters(1, 3)   # "One"
$r2 = $tr.Characters(5, 4)   # " Two"
$r3 = $tr.Characters(10, 6)  # " Three"

# These paragraphs carry no bullet/numbering (buNone).
$r1.ParagraphFormat.Bullet.Type = 0
$r2.ParagraphFormat.Bullet.Type = 0
$r3.ParagraphFormat.Bullet.Type = 0

# Apply the requested run formatting.
$r1.Font.Italic = $true
$r2.Font.Bold = $true
$r3.Font.Underline = $true
